$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Progress (D2): 1 -> 1.01
$ws.Range("D2").Value = 1.01

# DateAdded (E2): inline string "2025-07-21" -> real date serial 45859,
# formatted as a date. Apply the format twice (lowercase then uppercase)
# to match the two numFmt registrations (164/165) seen in the target file,
# with the cell ultimately styled with the uppercase "YYYY-MM-DD" format.
$ws.Range("E2").Value = 45859
$ws.Range("E2").NumberFormat = "yyyy-mm-dd"
$ws.Range("E2").NumberFormat = "YYYY-MM-DD"
